## Insert a new data row at row 189 (pushing existing rows 189-332 down to 190-333)
## and populate the new row with the record described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 189; this shifts all rows
# from 189 downward (189-332) one position down (becoming 190-333) and keeps
# the formatting of the following row for the new row's cells.
$ws.Rows.Item(189).Insert()

# Fill in the values for the newly inserted row 189.
$ws.Range("A189").Value = 8
$ws.Range("B189").Value = "Terminal La Palmera de La Serena"
$ws.Range("C189").Value = "Coquimbo"
$ws.Range("D189").Value = 44673
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 100114013
$ws.Range("G189").Value = "Zanahoria"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 600
$ws.Range("K189").Value = 6000
$ws.Range("L189").Value = 7000
$ws.Range("M189").Value = 6500
$ws.Range("N189").Value = '$/saco 20 kilos'
$ws.Range("O189").Value = "Provincia del Elquí"
$ws.Range("P189").Value = 325
$ws.Range("Q189").Value = 20
$ws.Range("R189").Value = "Hortaliza"
